$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 449 (shifts existing
# rows 449-461 down to 451-463) and fill them in with the new weekly
# "Crespo record" price observations dated 2021-09-09 (serial 44448).
$ws.Rows("449:450").Insert()

# New row 449: Crespo record / Primera -> Región Metropolitana
$ws.Cells.Item(449, 1).Value = 10
$ws.Cells.Item(449, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(449, 3).Value = "La Araucanía"
$ws.Cells.Item(449, 4).Value = 44448
$ws.Cells.Item(449, 5).Value = 9
$ws.Cells.Item(449, 6).Value = 100112006
$ws.Cells.Item(449, 7).Value = "Repollo"
$ws.Cells.Item(449, 8).Value = "Crespo record"
$ws.Cells.Item(449, 9).Value = "Primera"
$ws.Cells.Item(449, 10).Value = 4000
$ws.Cells.Item(449, 11).Value = 800
$ws.Cells.Item(449, 12).Value = 900
$ws.Cells.Item(449, 13).Value = 850
$ws.Cells.Item(449, 14).Value = "$/unidad"
$ws.Cells.Item(449, 15).Value = "Región Metropolitana"
$ws.Cells.Item(449, 16).Value = 850
$ws.Cells.Item(449, 17).Value = 1
$ws.Cells.Item(449, 18).Value = "Hortaliza"

# New row 450: Crespo record / Primera -> Región del Maule
$ws.Cells.Item(450, 1).Value = 10
$ws.Cells.Item(450, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(450, 3).Value = "La Araucanía"
$ws.Cells.Item(450, 4).Value = 44448
$ws.Cells.Item(450, 5).Value = 9
$ws.Cells.Item(450, 6).Value = 100112006
$ws.Cells.Item(450, 7).Value = "Repollo"
$ws.Cells.Item(450, 8).Value = "Crespo record"
$ws.Cells.Item(450, 9).Value = "Primera"
$ws.Cells.Item(450, 10).Value = 3500
$ws.Cells.Item(450, 11).Value = 800
$ws.Cells.Item(450, 12).Value = 900
$ws.Cells.Item(450, 13).Value = 843
$ws.Cells.Item(450, 14).Value = "$/unidad"
$ws.Cells.Item(450, 15).Value = "Región del Maule"
$ws.Cells.Item(450, 16).Value = 843
$ws.Cells.Item(450, 17).Value = 1
$ws.Cells.Item(450, 18).Value = "Hortaliza"
